$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 108 (shifts existing rows 108:189 down to 109:190,
# extending the used range to A1:R190), then populate it with the new
# weekly price record.
$ws.Rows.Item(108).Insert()

$ws.Range("A108").Value = 5
$ws.Range("B108").Value = "Macroferia Regional de Talca"
$ws.Range("C108").Value = "Maule"
$ws.Range("D108").Value = 44904
$ws.Range("E108").Value = 7
$ws.Range("F108").Value = 100112031
$ws.Range("G108").Value = "Poroto verde"
$ws.Range("H108").Value = "Sin especificar"
$ws.Range("I108").Value = "Primera"
$ws.Range("J108").Value = 150
$ws.Range("K108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("M108").Value = 30000
$ws.Range("N108").Value = "`$/saco 25 kilos"
$ws.Range("O108").Value = "Región del Maule"
$ws.Range("P108").Value = 1200
$ws.Range("Q108").Value = 25
$ws.Range("R108").Value = "Hortaliza"
